$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 = "本周工作" (this week's work) for the 8.7-8.11 column (C).
# Replace the draft text with the finalized write-up of completed work.
$ws.Range("C4").Value = "1. 后台及本地数据提取和分析；2. 嘉利、优闪订供货数据更新；3.启动和湖北乘用车供应商数据对接（路鹏）；4.参与询价单数据清洗、BI报表制作方案讨论；5.开始筹划建立平台车型库，并进行编码（初步计划对常见115乘用车品牌，14商用车品牌，以及1600余车系和16000余车型进行编码）"

# Row 6 = "下周计划" (next week's plan) for the 8.7-8.11 column (C), previously blank.
$ws.Range("C6").Value = "1. 基本数据分析；2.供货数据对接和更新（可能有优闪订、车享配、福耀和路鹏等）；3.完成对汽车品牌、车系编码，并细化车型定义和编码规则；4.与明觉、翱特交流数据配件数据结构等"

# Keep the on-screen selection in sync with where the edit happened.
$ws.Range("C4").Select()
